$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("2_")
$ws3 = $wb.Worksheets.Item("3_")
$ws7 = $wb.Worksheets.Item("7_Matching")

# --- Populate the "3_" sheet (new quiz content) ---
# Values are written in the same order the original author typed them so
# that new shared-string entries land at the expected indices.
$ws3.Range("A1").Value = 'We found the optimal launch angle in two ways: by manually sweeping through the "angle" parameter, and by using the SciPy function "minimize_scalare".   Choose which of those options has the advantage in each of the following categories.'
$ws3.Range("B1").Value = 'Correct order of definitions'
$ws3.Range("C1").Value = 'Definitions'

$ws3.Range("C2").Value = 'Manual sweeping is better for this'
$ws3.Range("C3").Value = 'A "black box" library function like "minimize_scalar" is better for this'

$ws3.Range("A2").Value = 'Accuracy of the answer'
$ws3.Range("A3").Value = 'Number of lines of code necessary to implement the method'
$ws3.Range("A4").Value = 'Speed of calculation'
$ws3.Range("A5").Value = 'Ease of fixing errors that arise'
$ws3.Range("A6").Value = 'Ease of understanding exactly what is happening in the code'

$ws3.Range("B2").Value = 'B'
$ws3.Range("B3").Value = 'B'
$ws3.Range("B4").Value = 'B'
$ws3.Range("B5").Value = 'A'
$ws3.Range("B6").Value = 'A'

$ws3.Range("D4").Value = 'A function like "minimize_scalar" is going to be optimized for speed by some really smart folks'
$ws3.Range("D3").Value = '"Minimize_scalar" is super efficient in terms of lines of code'
$ws3.Range("D5").Value = 'This might be a toss-up, but it can be really frustrating to figure out an unexpected error with a black box function'

# --- Formatting: copy the centered/wrapped style already used for column B
# elsewhere in the workbook so no new cellXfs entries are generated. ---
$ws7.Range("B1").Copy()
$ws3.Range("B1").PasteSpecial(-4122)
$ws7.Range("B2").Copy()
$ws3.Range("B2:B5").PasteSpecial(-4122)
$ws7.Range("B6").Copy()
$ws3.Range("B6:B7").PasteSpecial(-4122)

# --- Row heights (auto-computed heights from the original author's save) ---
$ws3.Rows.Item(1).RowHeight = 165
$ws3.Rows.Item(3).RowHeight = 45
$ws3.Rows.Item(4).RowHeight = 45
$ws3.Rows.Item(5).RowHeight = 60
$ws3.Rows.Item(6).RowHeight = 45

# --- Selections / active sheet ---
# "7_Matching" selection grows from A1:D7 to A1:D9 (its tab-selected state
# must not change, so this happens before we touch the active tab below).
$ws7.Range("A1:D9").Select()

# Move the old selection on "2_" to D8, and drop its tab-selected state by
# activating "3_" (which becomes the active tab, ending on a D6 selection).
$ws2.Range("D8").Select()
$ws3.Range("D6").Select()
